{"js": "// Cover Letter edits: recipient name/company/role/date/city changed, and the\n// two \"Amobee\" mentions + \"LinkedIn\" source renamed, per the commit diff.\n\nconst body = context.document.body;\n\n// Small helper: search for a unique (or Nth) occurrence of `text` and\n// replace it in place with `replacement`, preserving the run's formatting.\nasync function replaceOccurrence(searchText, replacement, occurrenceIndex = 0) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length <= occurrenceIndex) {\n    throw new Error(`Expected occurrence ${occurrenceIndex} of \"${searchText}\" not found`);\n  }\n\n  results.items[occurrenceIndex].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Recipient's first name: \"Olivia \" -> \"Andre \"\nawait replaceOccurrence(\"Olivia \", \"Andre \");\n\n// 2) Recipient's last name: \"Torrence\" -> \"Scheluchin\" (appears twice:\n//    once in the letterhead line, once in the \"Dear Ms./Mr.\" salutation).\nawait replaceOccurrence(\"Torrence\", \"Scheluchin\", 0);\nawait replaceOccurrence(\"Torrence\", \"Scheluchin\", 0); // list shrinks after each replace\n\n// 3) Letter date: \"June 12, 2020\" -> \"June 13, 2020\"\nawait replaceOccurrence(\"June 12, 2020\", \"June 13, 2020\");\n\n// 4) Recipient's title: \"HR Generalist\" -> \"Director of Talent at Insider Inc.\"\nawait replaceOccurrence(\"HR Generalist\", \"Director of Talent at Insider Inc.\");\n\n// 5) Recipient's city: \"Baltimore, Maryland\" -> \"New York, New York\"\nawait replaceOccurrence(\"Baltimore, Maryland\", \"New York, New York\");\n\n// 6) Salutation honorific: \"Dear Ms. \" -> \"Dear Mr. \"\nawait replaceOccurrence(\"Dear Ms. \", \"Dear Mr. \");\n\n// 7) Company name mentioned while introducing the application:\n//    \"...position at Amobee...\" -> \"...position at Insider Inc....\"\nawait replaceOccurrence(\"Amobee\", \"Insider Inc.\", 0);\n\n// 8) Job-posting source: \"LinkedIn\" -> \"Google Hire\"\nawait replaceOccurrence(\"LinkedIn\", \"Google Hire\");\n\n// 9) Company name mentioned again in the closing paragraph:\n//    \"...find at Amobee.\" -> \"...find at Insider.\"\n//    (the first \"Amobee\" is already gone, so this is now occurrence 0)\nawait replaceOccurrence(\"Amobee\", \"Insider\", 0);\n", "ps1": "# Cover Letter edits: recipient name/company/role/date/city changed, and the\n# two \"Amobee\" mentions + \"LinkedIn\" source renamed, per the commit diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# 1) Recipient's first name: \"Olivia \" -> \"Andre \"\nReplace-InRange $d.Content \"Olivia \" \"Andre \"\n\n# 2) Recipient's last name: \"Torrence\" -> \"Scheluchin\" (both the letterhead\n#    line and the \"Dear Ms./Mr.\" salutation use the same replacement).\nReplace-InRange $d.Content \"Torrence\" \"Scheluchin\"\n\n# 3) Letter date: \"June 12, 2020\" -> \"June 13, 2020\"\nReplace-InRange $d.Content \"June 12, 2020\" \"June 13, 2020\"\n\n# 4) Recipient's title: \"HR Generalist\" -> \"Director of Talent at Insider Inc.\"\nReplace-InRange $d.Content \"HR Generalist\" \"Director of Talent at Insider Inc.\"\n\n# 5) Recipient's city: \"Baltimore, Maryland\" -> \"New York, New York\"\nReplace-InRange $d.Content \"Baltimore, Maryland\" \"New York, New York\"\n\n# 6) Salutation honorific: \"Dear Ms. \" -> \"Dear Mr. \"\nReplace-InRange $d.Content \"Dear Ms. \" \"Dear Mr. \"\n\n# 7) Company name in the opening paragraph (6th paragraph): \"Amobee\" -> \"Insider Inc.\"\nReplace-InRange $d.Paragraphs.Item(6).Range \"Amobee\" \"Insider Inc.\"\n\n# 8) Job-posting source: \"LinkedIn\" -> \"Google Hire\"\nReplace-InRange $d.Content \"LinkedIn\" \"Google Hire\"\n\n# 9) Company name in the closing paragraph (12th paragraph): \"Amobee\" -> \"Insider\"\nReplace-InRange $d.Paragraphs.Item(12).Range \"Amobee\" \"Insider\"\n"}
